$wb = $excel.ActiveWorkbook

# Target OOXML column width is 17.2159881591797 characters-units, but the
# host's ColumnWidth setter quantizes to 1/6-character steps, so feed it
# the character-width value whose quantized result lands closest to the
# target stored width (17.166666666666668 vs. 17.2159881591797).
$colWidth = 16.3333333333333

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 07:03:08"
$wsOverview.Columns.Item(5).ColumnWidth = $colWidth
$wsOverview.Columns.Item(6).ColumnWidth = $colWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-30 07:02:57"
$wsZhCn.Columns.Item(3).ColumnWidth = $colWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-30 07:03:08"
$wsDeDe.Columns.Item(3).ColumnWidth = $colWidth
